$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update time_taken timestamps in the "data" sheet (F2:F46)
$ws.Range("F2").Value = "2021-10-05 14:35:07.336765"
$ws.Range("F3").Value = "2021-10-05 14:35:07.336773"
$ws.Range("F4").Value = "2021-10-05 14:35:07.336776"
$ws.Range("F5").Value = "2021-10-05 14:35:07.336779"
$ws.Range("F6").Value = "2021-10-05 14:35:07.336782"
$ws.Range("F7").Value = "2021-10-05 14:35:07.336785"
$ws.Range("F8").Value = "2021-10-05 14:35:07.336788"
$ws.Range("F9").Value = "2021-10-05 14:35:07.336790"
$ws.Range("F10").Value = "2021-10-05 14:35:07.336793"
$ws.Range("F11").Value = "2021-10-05 14:35:07.336796"
$ws.Range("F12").Value = "2021-10-05 14:35:07.336798"
$ws.Range("F13").Value = "2021-10-05 14:35:07.336801"
$ws.Range("F14").Value = "2021-10-05 14:35:07.336804"
$ws.Range("F15").Value = "2021-10-05 14:35:07.336806"
$ws.Range("F16").Value = "2021-10-05 14:35:07.336809"
$ws.Range("F17").Value = "2021-10-05 14:35:07.336811"
$ws.Range("F18").Value = "2021-10-05 14:35:07.336814"
$ws.Range("F19").Value = "2021-10-05 14:35:07.336817"
$ws.Range("F20").Value = "2021-10-05 14:35:07.336820"
$ws.Range("F21").Value = "2021-10-05 14:35:07.336822"
$ws.Range("F22").Value = "2021-10-05 14:35:07.336825"
$ws.Range("F23").Value = "2021-10-05 14:35:07.336827"
$ws.Range("F24").Value = "2021-10-05 14:35:07.336830"
$ws.Range("F25").Value = "2021-10-05 14:35:07.336832"
$ws.Range("F26").Value = "2021-10-05 14:35:07.336835"
$ws.Range("F27").Value = "2021-10-05 14:35:07.336838"
$ws.Range("F28").Value = "2021-10-05 14:35:07.336841"
$ws.Range("F29").Value = "2021-10-05 14:35:07.336843"
$ws.Range("F30").Value = "2021-10-05 14:35:07.336846"
$ws.Range("F31").Value = "2021-10-05 14:35:07.336848"
$ws.Range("F32").Value = "2021-10-05 14:35:07.336851"
$ws.Range("F33").Value = "2021-10-05 14:35:07.336854"
$ws.Range("F34").Value = "2021-10-05 14:35:07.336857"
$ws.Range("F35").Value = "2021-10-05 14:35:07.336859"
$ws.Range("F36").Value = "2021-10-05 14:35:07.336862"
$ws.Range("F37").Value = "2021-10-05 14:35:07.336865"
$ws.Range("F38").Value = "2021-10-05 14:35:07.336867"
$ws.Range("F39").Value = "2021-10-05 14:35:07.336870"
$ws.Range("F40").Value = "2021-10-05 14:35:07.336872"
$ws.Range("F41").Value = "2021-10-05 14:35:07.336875"
$ws.Range("F42").Value = "2021-10-05 14:35:07.336878"
$ws.Range("F43").Value = "2021-10-05 14:35:07.336881"
$ws.Range("F44").Value = "2021-10-05 14:35:07.336883"
$ws.Range("F45").Value = "2021-10-05 14:35:07.336886"
$ws.Range("F46").Value = "2021-10-05 14:35:07.336888"

# Add the new "metadata" worksheet after "data"
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws)
$ws2.Name = "metadata"

# Copy the header style (bold, centered, bordered) used on sheet "data" row 1
$ws.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

# Copy the index-column style (bold, centered, bordered) used on sheet "data" column A
$ws.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Optic Atrophy"
$ws2.Range("C2").Value = 149
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.0"
$ws.Range("B2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)
$ws2.Range("E2").Value = "2021-07-14T08:39:04.930772Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:07.332904"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/149/?format=json"

[void]$ws.Activate()
[void]$ws.Range("A1").Select()
